{"js": "// Each entry is [row, col, oldText, newText] for the single 5-column table\n// (0-based row/col, matching Office.js's table.getCell indexing). The blank\n// spacer rows (1-3, 5-7, 9-11, ...) are simply not listed. Cross-checking\n// the old text before writing guards against addressing the wrong cell.\nconst cellEdits = [\n  [0, 0, \"89\u00f73=\", \"63\u00f77=\"],\n  [0, 1, \"23\u00f79=\", \"63\u00f79=\"],\n  [0, 2, \"70\u00f72=\", \"31\u00f77=\"],\n  [0, 3, \"47\u00f72=\", \"20\u00f75=\"],\n  [0, 4, \"71\u00f72=\", \"82\u00f76=\"],\n\n  [4, 0, \"71\u00f79=\", \"11\u00f79=\"],\n  [4, 1, \"90\u00f72=\", \"24\u00f72=\"],\n  [4, 2, \"49\u00f77=\", \"36\u00f78=\"],\n  [4, 3, \"89\u00f78=\", \"78\u00f75=\"],\n  [4, 4, \"37\u00f78=\", \"64\u00f73=\"],\n\n  [8, 0, \"42\u00f74=\", \"78\u00f78=\"],\n  [8, 1, \"73\u00f77=\", \"19\u00f76=\"],\n  [8, 2, \"18\u00f74=\", \"65\u00f79=\"],\n  [8, 3, \"79\u00f73=\", \"20\u00f77=\"],\n  [8, 4, \"53\u00f78=\", \"17\u00f77=\"],\n\n  [12, 0, \"28\u00f77=\", \"17\u00f73=\"],\n  [12, 1, \"45\u00f72=\", \"64\u00f72=\"],\n  [12, 2, \"82\u00f76=\", \"33\u00f75=\"],\n  [12, 3, \"42\u00f77=\", \"61\u00f72=\"],\n  [12, 4, \"72\u00f72=\", \"50\u00f73=\"],\n\n  [16, 0, \"94\u00f75=\", \"39\u00f75=\"],\n  [16, 1, \"27\u00f72=\", \"32\u00f74=\"],\n  [16, 2, \"81\u00f76=\", \"78\u00f76=\"],\n  [16, 3, \"11\u00f73=\", \"47\u00f73=\"],\n  [16, 4, \"85\u00f76=\", \"77\u00f75=\"],\n];\n\n// 1) Update the date line (first paragraph of the body).\nconst dateRange = context.document.body.paragraphs.getFirst().getRange();\ndateRange.load(\"text\");\n\n// 2) Gather the range for every table cell to edit, addressed by\n//    row/column position (not by searching for its text) so that values\n//    which coincide with another cell's *new* text are never\n//    double-replaced.\nconst table = context.document.body.tables.getFirst();\nconst editRanges = cellEdits.map(([row, col, oldText, newText]) => {\n  const range = table.getCell(row, col).body.paragraphs.getFirst().getRange();\n  range.load(\"text\");\n  return { range, oldText, newText };\n});\n\nawait context.sync();\n\nif (dateRange.text === \"2025-07-13 Sunday\") {\n  dateRange.insertText(\"2025-07-14 Monday\", Word.InsertLocation.replace);\n}\n\nfor (const { range, oldText, newText } of editRanges) {\n  if (range.text === oldText) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every division-problem cell in the single\n# table, addressed by position (Paragraphs(1) / Table(1).Cell(row,col))\n# rather than by searching for old text, so that values which coincide\n# with another cell's *new* text are never double-replaced. Each write\n# is guarded by a check against the expected old value.\n\n$d = $word.ActiveDocument\n\n# 1) Date line - first paragraph in the document body. Range.Text carries a\n#    trailing paragraph mark (CR, char 13), so trim it before comparing.\n$dateText = $d.Paragraphs(1).Range.Text.TrimEnd([char]13, [char]7)\nif ($dateText -eq \"2025-07-13 Sunday\") {\n    $d.Paragraphs(1).Range.Text = \"2025-07-14 Monday\"\n}\n\n# 2) Division problems - single 5-column table; every 4th row (1,5,9,13,17)\n#    holds the visible problems, the rows between are blank spacers.\n$tbl = $d.Tables(1)\n\n$cellEdits = @(\n    @{ Row = 1; Col = 1; Old = \"89\u00f73=\"; New = \"63\u00f77=\" }\n    @{ Row = 1; Col = 2; Old = \"23\u00f79=\"; New = \"63\u00f79=\" }\n    @{ Row = 1; Col = 3; Old = \"70\u00f72=\"; New = \"31\u00f77=\" }\n    @{ Row = 1; Col = 4; Old = \"47\u00f72=\"; New = \"20\u00f75=\" }\n    @{ Row = 1; Col = 5; Old = \"71\u00f72=\"; New = \"82\u00f76=\" }\n\n    @{ Row = 5; Col = 1; Old = \"71\u00f79=\"; New = \"11\u00f79=\" }\n    @{ Row = 5; Col = 2; Old = \"90\u00f72=\"; New = \"24\u00f72=\" }\n    @{ Row = 5; Col = 3; Old = \"49\u00f77=\"; New = \"36\u00f78=\" }\n    @{ Row = 5; Col = 4; Old = \"89\u00f78=\"; New = \"78\u00f75=\" }\n    @{ Row = 5; Col = 5; Old = \"37\u00f78=\"; New = \"64\u00f73=\" }\n\n    @{ Row = 9; Col = 1; Old = \"42\u00f74=\"; New = \"78\u00f78=\" }\n    @{ Row = 9; Col = 2; Old = \"73\u00f77=\"; New = \"19\u00f76=\" }\n    @{ Row = 9; Col = 3; Old = \"18\u00f74=\"; New = \"65\u00f79=\" }\n    @{ Row = 9; Col = 4; Old = \"79\u00f73=\"; New = \"20\u00f77=\" }\n    @{ Row = 9; Col = 5; Old = \"53\u00f78=\"; New = \"17\u00f77=\" }\n\n    @{ Row = 13; Col = 1; Old = \"28\u00f77=\"; New = \"17\u00f73=\" }\n    @{ Row = 13; Col = 2; Old = \"45\u00f72=\"; New = \"64\u00f72=\" }\n    @{ Row = 13; Col = 3; Old = \"82\u00f76=\"; New = \"33\u00f75=\" }\n    @{ Row = 13; Col = 4; Old = \"42\u00f77=\"; New = \"61\u00f72=\" }\n    @{ Row = 13; Col = 5; Old = \"72\u00f72=\"; New = \"50\u00f73=\" }\n\n    @{ Row = 17; Col = 1; Old = \"94\u00f75=\"; New = \"39\u00f75=\" }\n    @{ Row = 17; Col = 2; Old = \"27\u00f72=\"; New = \"32\u00f74=\" }\n    @{ Row = 17; Col = 3; Old = \"81\u00f76=\"; New = \"78\u00f76=\" }\n    @{ Row = 17; Col = 4; Old = \"11\u00f73=\"; New = \"47\u00f73=\" }\n    @{ Row = 17; Col = 5; Old = \"85\u00f76=\"; New = \"77\u00f75=\" }\n)\n\nforeach ($edit in $cellEdits) {\n    $cellRange = $tbl.Cell($edit.Row, $edit.Col).Range\n    # A table cell's Range.Text carries trailing CR + cell-mark (char 13, 7).\n    $cellText = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($cellText -eq $edit.Old) {\n        $cellRange.Text = $edit.New\n    }\n}\n"}
